{"js": "// Delete the first paragraph of the document body\n// (\"Video provides a powerful way to help you prove your point. ...\")\n// including its paragraph mark, leaving the remaining paragraphs intact.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst firstParagraph = paragraphs.items[0];\nfirstParagraph.load(\"text\");\nawait context.sync();\n\nif (firstParagraph.text.indexOf(\"Video provides a powerful way\") !== -1) {\n  firstParagraph.delete();\n  await context.sync();\n}\n", "ps1": "# Delete the first paragraph of the document body\n# (\"Video provides a powerful way to help you prove your point. ...\")\n# including its paragraph mark, leaving the remaining paragraphs intact.\n\n$d = $word.ActiveDocument\n$firstParagraph = $d.Paragraphs.Item(1)\n\nif ($firstParagraph.Range.Text -like \"*Video provides a powerful way*\") {\n    $firstParagraph.Range.Delete()\n}\n"}
